$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '43.779.03'
$ws.Range('E2').Value = '  -0.12%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.291.40'
$ws.Range('E3').Value = '  -1.16%  '
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '103.43'
$ws.Range('E5').Value = '  +5.91%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '270.54'
$ws.Range('E6').Value = '  -0.78%  '
$ws.Range('E7').Value = '  -1.60%  '
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.609'
$ws.Range('E9').Value = '  -2.80%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '45.93'
$ws.Range('E10').Value = '  +1.01%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0936'
$ws.Range('E11').Value = '  -1.50%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '8.15'
$ws.Range('E12').Value = '  +1.44%  '
$ws.Range('E13').Value = '  +1.43%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '15.57'
$ws.Range('E14').Value = '  +0.38%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.850'
$ws.Range('E15').Value = '  -2.99%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.291.55'
$ws.Range('E16').Value = '  -1.24%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '43.761.03'
$ws.Range('E17').Value = '  -0.03%  '
$ws.Range('E18').Value = '  +0.41%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.28'
$ws.Range('E19').Value = '  -2.07%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '72.32'
$ws.Range('E20').Value = '  -1.44%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '2.52'
$ws.Range('E21').Value = '  +10.75%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '232.27'
$ws.Range('E22').Value = '  -3.28%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.87'
$ws.Range('E23').Value = '  +12.65%  '
$ws.Range('E24').Value = '  -1.77%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '11.30'
$ws.Range('E26').Value = '  -0.79%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '40.85'
$ws.Range('E27').Value = '  +6.81%  '
$ws.Range('E28').Value = '  -1.58%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.22'
$ws.Range('E29').Value = '  -2.92%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '177.76'
$ws.Range('E30').Value = '  +1.45%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '21.85'
$ws.Range('E31').Value = '  -2.67%  '
$ws.Range('E32').Value = '  -1.98%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.53'
$ws.Range('E33').Value = '  +0.70%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.90'
$ws.Range('E34').Value = '  +10.63%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.126'
$ws.Range('E35').Value = '  -1.09%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.110'
$ws.Range('E36').Value = '  +0.22%  '
$ws.Range('E37').Value = '  -0.56%  '
$ws.Range('E38').Value = '  +4.00%  '
$ws.Range('E39').Value = '  -3.46%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.34'
$ws.Range('E40').Value = '  -2.48%  '
$ws.Range('E41').Value = '  -2.91%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '65.52'
$ws.Range('E42').Value = '  +4.26%  '
$ws.Range('E43').Value = '  -0.16%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.30'
$ws.Range('E44').Value = '  -1.21%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '8.82'
$ws.Range('E45').Value = '  -4.43%  '
$ws.Range('E46').Value = '  -1.06%  '
$ws.Range('E47').Value = '  +2.48%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '99.24'
$ws.Range('E48').Value = '  -1.28%  '
$ws.Range('E49').Value = '  +10.92%  '
$ws.Range('E50').Value = '  +4.83%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.524.85'
$ws.Range('E51').Value = '  -0.79%  '
